# Fruta / hortaliza, semanal
# Two new weekly price records (for the most recent sampling date) are
# added to the top of the data block for this market/product, pushing the
# existing rows down by two positions (row 120 -> 122, ... row 183 -> 185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 120/121; everything from the old row 120 down
# shifts to row 122 onward (Excel's normal "insert row" semantics, which
# also grows the sheet's used range / <dimension> automatically).
$ws.Range("A120:A121").EntireRow.Insert()

# New row 120 - "Primera" quality, newest sampling date (2023-06-20)
$ws.Cells.Item(120, 1).Value  = 1
$ws.Cells.Item(120, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(120, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(120, 4).Value  = 45097
$ws.Cells.Item(120, 5).Value  = 15
$ws.Cells.Item(120, 6).Value  = 'Fruta'
$ws.Cells.Item(120, 7).Value  = 100108
$ws.Cells.Item(120, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(120, 9).Value  = 100108003
$ws.Cells.Item(120, 10).Value = 'Maracuyá'
$ws.Cells.Item(120, 11).Value = 'Sin especificar'
$ws.Cells.Item(120, 12).Value = 'Primera'
$ws.Cells.Item(120, 13).Value = 100
$ws.Cells.Item(120, 14).Value = 27000
$ws.Cells.Item(120, 15).Value = 28000
$ws.Cells.Item(120, 16).Value = 27500
$ws.Cells.Item(120, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(120, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(120, 19).Value = 1375
$ws.Cells.Item(120, 20).Value = 20

# New row 121 - "Segunda" quality, same sampling date
$ws.Cells.Item(121, 1).Value  = 1
$ws.Cells.Item(121, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(121, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(121, 4).Value  = 45097
$ws.Cells.Item(121, 5).Value  = 15
$ws.Cells.Item(121, 6).Value  = 'Fruta'
$ws.Cells.Item(121, 7).Value  = 100108
$ws.Cells.Item(121, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(121, 9).Value  = 100108003
$ws.Cells.Item(121, 10).Value = 'Maracuyá'
$ws.Cells.Item(121, 11).Value = 'Sin especificar'
$ws.Cells.Item(121, 12).Value = 'Segunda'
$ws.Cells.Item(121, 13).Value = 120
$ws.Cells.Item(121, 14).Value = 23000
$ws.Cells.Item(121, 15).Value = 24000
$ws.Cells.Item(121, 16).Value = 23500
$ws.Cells.Item(121, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(121, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(121, 19).Value = 1175
$ws.Cells.Item(121, 20).Value = 20
